$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 44138
$ws.Range("D2").Value = 66299658
$ws.Range("C3").Value = 104366
$ws.Range("D3").Value = 159539928
$ws.Range("C4").Value = 35606
$ws.Range("D4").Value = 56672007
$ws.Range("C5").Value = 10273
$ws.Range("D5").Value = 17012107
$ws.Range("C6").Value = 2680
$ws.Range("D6").Value = 4693015
$ws.Range("C7").Value = 304
$ws.Range("D7").Value = 583683
$ws.Range("C12").Value = 46772
$ws.Range("D12").Value = 63969879
$ws.Range("C13").Value = 11151
$ws.Range("D13").Value = 16493467
$ws.Range("C14").Value = 29171
$ws.Range("D14").Value = 43569108
$ws.Range("C15").Value = 9249
$ws.Range("D15").Value = 14128875
$ws.Range("C16").Value = 2469
$ws.Range("D16").Value = 3869199
$ws.Range("C17").Value = 536
$ws.Range("D17").Value = 842282
$ws.Range("C20").Value = 11475
$ws.Range("D20").Value = 15184768
$ws.Range("C21").Value = 15171
$ws.Range("D21").Value = 22069321
$ws.Range("C22").Value = 35248
$ws.Range("D22").Value = 52089859
$ws.Range("C23").Value = 11363
$ws.Range("D23").Value = 17234711
$ws.Range("C24").Value = 3039
$ws.Range("D24").Value = 4704504
$ws.Range("C25").Value = 674
$ws.Range("D25").Value = 1066465
$ws.Range("C26").Value = 59
$ws.Range("D26").Value = 124627
$ws.Range("C27").Value = 12951
$ws.Range("D27").Value = 17160576
$ws.Range("C28").Value = 8935
$ws.Range("D28").Value = 13226971
$ws.Range("C29").Value = 25594
$ws.Range("D29").Value = 38351312
$ws.Range("C30").Value = 8756
$ws.Range("D30").Value = 13377026
$ws.Range("C31").Value = 2240
$ws.Range("D31").Value = 3427863
$ws.Range("C32").Value = 483
$ws.Range("D32").Value = 759258
$ws.Range("C34").Value = 9374
$ws.Range("D34").Value = 12364637
$ws.Range("C35").Value = 3888
$ws.Range("D35").Value = 5790804
$ws.Range("C36").Value = 9069
$ws.Range("D36").Value = 13737351
$ws.Range("C37").Value = 3556
$ws.Range("D37").Value = 5477403
$ws.Range("C38").Value = 907
$ws.Range("D38").Value = 1385216
$ws.Range("C39").Value = 196
$ws.Range("D39").Value = 325686
$ws.Range("C41").Value = 2906
$ws.Range("D41").Value = 3938931
$ws.Range("C42").Value = 19833
$ws.Range("D42").Value = 29337671
$ws.Range("C43").Value = 57437
$ws.Range("D43").Value = 85749057
$ws.Range("C44").Value = 20987
$ws.Range("D44").Value = 31744852
$ws.Range("C45").Value = 6328
$ws.Range("D45").Value = 9610595
$ws.Range("C46").Value = 1566
$ws.Range("D46").Value = 2465768
$ws.Range("C47").Value = 104
$ws.Range("D47").Value = 208961
$ws.Range("C50").Value = 19149
$ws.Range("D50").Value = 25467921
$ws.Range("C51").Value = 2475
$ws.Range("D51").Value = 3700140
$ws.Range("C52").Value = 8287
$ws.Range("D52").Value = 12512103
$ws.Range("C53").Value = 2785
$ws.Range("D53").Value = 4390839
$ws.Range("C54").Value = 883
$ws.Range("D54").Value = 1386723
$ws.Range("C55").Value = 252
$ws.Range("D55").Value = 440110
$ws.Range("C56").Value = 28
$ws.Range("D56").Value = 76000
$ws.Range("C57").Value = 8184
$ws.Range("D57").Value = 11332319
$ws.Range("C58").Value = 1711
$ws.Range("D58").Value = 3535516
$ws.Range("C59").Value = 4067
$ws.Range("D59").Value = 8333064
$ws.Range("C60").Value = 1621
$ws.Range("D60").Value = 3356997
$ws.Range("C61").Value = 540
$ws.Range("D61").Value = 1097843
$ws.Range("C64").Value = 2656
$ws.Range("D64").Value = 5056445
$ws.Range("C65").Value = 18080
$ws.Range("D65").Value = 27082856
$ws.Range("C66").Value = 51136
$ws.Range("D66").Value = 77816560
$ws.Range("C67").Value = 17804
$ws.Range("D67").Value = 27908095
$ws.Range("C68").Value = 5237
$ws.Range("D68").Value = 8344594
$ws.Range("C69").Value = 1213
$ws.Range("D69").Value = 2119736
$ws.Range("C70").Value = 126
$ws.Range("D70").Value = 257582
$ws.Range("C71").Value = 20
$ws.Range("D71").Value = 26285
$ws.Range("C73").Value = 16949
$ws.Range("D73").Value = 22312883
$ws.Range("C74").Value = 68497
$ws.Range("D74").Value = 108582181
$ws.Range("C75").Value = 183958
$ws.Range("D75").Value = 294798185
$ws.Range("C76").Value = 78488
$ws.Range("D76").Value = 131798798
$ws.Range("C77").Value = 25991
$ws.Range("D77").Value = 45686531
$ws.Range("C78").Value = 7349
$ws.Range("D78").Value = 14967815
$ws.Range("C79").Value = 625
$ws.Range("D79").Value = 1702113
$ws.Range("C85").Value = 67117
$ws.Range("D85").Value = 93443511
$ws.Range("C86").Value = 5318
$ws.Range("D86").Value = 7787931
$ws.Range("C87").Value = 12971
$ws.Range("D87").Value = 19356855
$ws.Range("C88").Value = 4219
$ws.Range("D88").Value = 6388098
$ws.Range("C89").Value = 1490
$ws.Range("D89").Value = 2226111
$ws.Range("C90").Value = 367
$ws.Range("D90").Value = 572012
$ws.Range("C93").Value = 6045
$ws.Range("D93").Value = 8133170
$ws.Range("C94").Value = 1932
$ws.Range("D94").Value = 2843772
$ws.Range("C95").Value = 6139
$ws.Range("D95").Value = 9303867
$ws.Range("C96").Value = 2182
$ws.Range("D96").Value = 3349191
$ws.Range("C97").Value = 799
$ws.Range("D97").Value = 1219797
$ws.Range("C98").Value = 251
$ws.Range("D98").Value = 409758
$ws.Range("C101").Value = 4165
$ws.Range("D101").Value = 5589115
$ws.Range("C102").Value = 951
$ws.Range("D102").Value = 1869288
$ws.Range("C103").Value = 655
$ws.Range("D103").Value = 1385476
$ws.Range("C107").Value = 7
$ws.Range("D107").Value = 25485
$ws.Range("C108").Value = 12704
$ws.Range("D108").Value = 19138267
$ws.Range("C109").Value = 32633
$ws.Range("D109").Value = 49254575
$ws.Range("C110").Value = 10929
$ws.Range("D110").Value = 16755717
$ws.Range("C111").Value = 3105
$ws.Range("D111").Value = 4847624
$ws.Range("C112").Value = 647
$ws.Range("D112").Value = 1052004
$ws.Range("C113").Value = 100
$ws.Range("D113").Value = 248883
$ws.Range("C116").Value = 10906
$ws.Range("D116").Value = 14419538
$ws.Range("C117").Value = 35205
$ws.Range("D117").Value = 51925686
$ws.Range("C118").Value = 74268
$ws.Range("D118").Value = 110467030
$ws.Range("C119").Value = 23791
$ws.Range("D119").Value = 36030340
$ws.Range("C120").Value = 6808
$ws.Range("D120").Value = 10413600
$ws.Range("C121").Value = 1446
$ws.Range("D121").Value = 2345992
$ws.Range("C122").Value = 151
$ws.Range("D122").Value = 267291
$ws.Range("C126").Value = 28728
$ws.Range("D126").Value = 38414891
$ws.Range("C127").Value = 42349
$ws.Range("D127").Value = 63472603
$ws.Range("C128").Value = 87821
$ws.Range("D128").Value = 133490122
$ws.Range("C129").Value = 26882
$ws.Range("D129").Value = 42628108
$ws.Range("C130").Value = 7376
$ws.Range("D130").Value = 11802750
$ws.Range("C131").Value = 1666
$ws.Range("D131").Value = 2920314
$ws.Range("C135").Value = 35737
$ws.Range("D135").Value = 47704161
$ws.Range("C136").Value = 15429
$ws.Range("D136").Value = 22644124
$ws.Range("C137").Value = 36460
$ws.Range("D137").Value = 53951965
$ws.Range("C138").Value = 12875
$ws.Range("D138").Value = 19274381
$ws.Range("C139").Value = 3422
$ws.Range("D139").Value = 5249327
$ws.Range("C140").Value = 661
$ws.Range("D140").Value = 1093808
$ws.Range("C141").Value = 61
$ws.Range("D141").Value = 128655
$ws.Range("C144").Value = 12113
$ws.Range("D144").Value = 16223303
$ws.Range("C145").Value = 41941
$ws.Range("D145").Value = 64069767
$ws.Range("C146").Value = 95789
$ws.Range("D146").Value = 149091499
$ws.Range("C147").Value = 28658
$ws.Range("D147").Value = 46222289
$ws.Range("C148").Value = 7767
$ws.Range("D148").Value = 13119022
$ws.Range("C149").Value = 2013
$ws.Range("D149").Value = 3670859
$ws.Range("C150").Value = 182
$ws.Range("D150").Value = 418276
$ws.Range("C151").Value = 24
$ws.Range("D151").Value = 61500
$ws.Range("C152").Value = 33346
$ws.Range("D152").Value = 45401686
